$wb = $excel.ActiveWorkbook

# Sheet2 ("Sheet2") has a list of products. Replace "Apple iphone" with "Sunglasses".
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A4").Value = "Sunglasses"
